# Refresh the crypto price/volume columns (D = Price, E = Volume(1h))
# to the latest scrape. Price strings that look numeric ("219.10",
# "11.00", ...) are written with a leading apostrophe so Excel keeps
# them as literal text (matching the source data, which stores these
# as plain strings, e.g. multi-dot "26.591.23" would not even parse as
# a number) instead of coercing to a number and dropping formatting
# such as trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.591.23'
$ws.Range('E2').Value = '  -7.22%  '
$ws.Range('D3').Value = '1.695.23'
$ws.Range('E3').Value = '  -5.84%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = '''219.10'
$ws.Range('E5').Value = '  -5.38%  '
$ws.Range('D6').Value = '''0.5089'
$ws.Range('E6').Value = '  -13.69%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').Value = '''0.2645'
$ws.Range('E8').Value = '  -4.48%  '
$ws.Range('E9').Value = '  -4.91%  '
$ws.Range('D10').Value = '''0.06290'
$ws.Range('D11').Value = '''0.07374'
$ws.Range('E11').Value = '  -1.86%  '
$ws.Range('D12').Value = '1.697.02'
$ws.Range('E12').Value = '  -7.59%  '
$ws.Range('D13').Value = '''4.512'
$ws.Range('E13').Value = '  -5.28%  '
$ws.Range('D14').Value = '''0.5829'
$ws.Range('D15').Value = '1.926.38'
$ws.Range('E15').Value = '  -5.78%  '
$ws.Range('D16').Value = '''0.000008388'
$ws.Range('E16').Value = '  -8.14%  '
$ws.Range('D17').Value = '''65.49'
$ws.Range('E17').Value = '  -13.29%  '
$ws.Range('D18').Value = '26.619.92'
$ws.Range('E18').Value = '  -7.04%  '
$ws.Range('E19').Value = '  -8.20%  '
$ws.Range('D20').Value = '''1.005'
$ws.Range('E20').Value = '  +0.23%  '
$ws.Range('D21').Value = '''11.00'
$ws.Range('E21').Value = '  -4.24%  '
$ws.Range('D22').Value = '''186.31'
$ws.Range('E22').Value = '  -11.34%  '
$ws.Range('D23').Value = '''6.261'
$ws.Range('E23').Value = '  -8.16%  '
$ws.Range('E24').Value = '  +0.27%  '
$ws.Range('D25').Value = '''144.46'
$ws.Range('E25').Value = '  -6.01%  '
$ws.Range('D26').Value = '''7.516'
$ws.Range('E26').Value = '  -4.24%  '
$ws.Range('D27').Value = '''0.1157'
$ws.Range('E27').Value = '  -8.64%  '
$ws.Range('E28').Value = '  -4.54%  '
$ws.Range('D29').Value = '''1.343'
$ws.Range('E29').Value = '  -5.13%  '
$ws.Range('D30').Value = '''0.05655'
$ws.Range('E30').Value = '  -8.42%  '
$ws.Range('D31').Value = '''1.336'
$ws.Range('E31').Value = '  -6.04%  '
$ws.Range('D32').Value = '''3.511'
$ws.Range('E32').Value = '  -6.93%  '
$ws.Range('D33').Value = '''3.481'
$ws.Range('E33').Value = '  -7.77%  '
$ws.Range('D34').Value = '''1.640'
$ws.Range('E34').Value = '  -5.17%  '
$ws.Range('D35').Value = '''1.020'
$ws.Range('E35').Value = '  -3.54%  '
$ws.Range('D36').Value = '''0.6034'
$ws.Range('E36').Value = '  -5.90%  '
$ws.Range('D37').Value = '''2.364'
$ws.Range('E37').Value = '  -5.42%  '
$ws.Range('D38').Value = '''2.681'
$ws.Range('E38').Value = '  -1.06%  '
$ws.Range('E39').Value = '  -4.90%  '
$ws.Range('D40').Value = '1.094.62'
$ws.Range('E40').Value = '  -4.52%  '
$ws.Range('D41').Value = '''0.8584'
$ws.Range('E41').Value = '  -2.69%  '
$ws.Range('E42').Value = '  -10.63%  '
$ws.Range('D43').Value = '''1.002'
$ws.Range('E43').Value = '  -0.46%  '
$ws.Range('D44').Value = '''99.60'
$ws.Range('E44').Value = '  -0.49%  '
$ws.Range('D45').Value = '1.852.69'
$ws.Range('E45').Value = '  -5.15%  '
$ws.Range('D46').Value = '''0.00000000109'
$ws.Range('E46').Value = '  -1.35%  '
$ws.Range('D47').Value = '''56.70'
$ws.Range('E47').Value = '  -5.92%  '
$ws.Range('D48').Value = '''8.198'
$ws.Range('E48').Value = '  -1.91%  '
$ws.Range('D49').Value = '''1.005'
$ws.Range('E49').Value = '  +0.29%  '
$ws.Range('D50').Value = '''0.05241'
$ws.Range('E50').Value = '  -4.02%  '
$ws.Range('D51').Value = '''0.4325'
$ws.Range('E51').Value = '  -3.50%  '
